$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected cell/view on the worksheet
$ws.Range("I22").Select()

# Update the Target values for rows 14-25 (column B)
$ws.Range("B14").Value = 459193710.18150002
$ws.Range("B15").Value = 422540529.20887506
$ws.Range("B16").Value = 137691163.72475001
$ws.Range("B17").Value = 552065941.25979996
$ws.Range("B18").Value = 527964562.28925002
$ws.Range("B19").Value = 545770222.41162503
$ws.Range("B20").Value = 428309905.52437502
$ws.Range("B21").Value = 592410547.87699997
$ws.Range("B22").Value = 445800000
$ws.Range("B23").Value = 514956287.90000004
$ws.Range("B24").Value = 508654001.75
$ws.Range("B25").Value = 548376614.64999998

# Update window size/position of the workbook view
$win = $excel.ActiveWindow
$win.Left = 1950
$win.Top = 1200
$win.Width = 27855
$win.Height = 31200
